# Update the "Förändrad" date column (C) from 2023-10-06 (45205) to 2023-10-07 (45206)
# for all data rows (2 through 39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
